# "1st changes of mifos to finflux"
#
# On the "Repayment schedule" sheet, insert a new blank column before
# column N (the existing N/O/P columns - "Late", "heading", "Outstanding" -
# shift right to O/P/Q). The new column inherits the column width of the
# column immediately to its left (M), matching Excel's default
# insert-column behaviour.
#
# Finally, switch the active sheet/selection from "Transactions" to
# "Repayment schedule".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

$ws.Activate()
$ws.Range("R7").Select()
